# Se agrega nueva mercancia e inventario
# Update the inventory/merchandise codes on the active sheet ("Semilla 11").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the three merchandise/inventory codes in column B (rows 9-11).
$ws.Range("B9").Value = "708603165"
$ws.Range("B10").Value = "325869013"
$ws.Range("B11").Value = "702923575"

# Leave the selection on the last touched cell, matching the saved view state.
$ws.Range("B12").Select()
